# Update the "as_of_utc" timestamp column (AA) for rows 2-26 on both the
# "Главные" (index 2) and "Линейные" (index 3) worksheets, changing the
# stamp from 2025-11-04 09:23:57 to 2025-11-04 09:58:43.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-04 09:58:43"

$sheetIndexes = @(2, 3)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
